# Updated cryptos list - apply latest price/volume data scrape

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns store plain text (e.g. "0.7120", "  +0.58%  ").
# Force text format before writing so Excel doesn't auto-coerce numeric-
# looking strings into real numbers, then restore the original (unstyled)
# look so no stray number formatting lingers on the cells.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "29.308.45"
$ws.Range("E2").Value = "  +0.55%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.871.92"
$ws.Range("E3").Value = "  +0.61%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 - XRP
$ws.Range("D5").Value = "0.7116"
$ws.Range("E5").Value = "  +0.76%  "

# Row 6 - BNB
$ws.Range("D6").Value = "241.63"
$ws.Range("E6").Value = "  +0.30%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.06%  "

# Row 8 - was Dogecoin, now Cardano (rows 8 and 9 swapped content)
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "0.3109"
$ws.Range("E8").Value = "  +0.91%  "

# Row 9 - was Cardano, now Dogecoin
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.07780"
$ws.Range("E9").Value = "  +2.10%  "

# Row 10 - Solana
$ws.Range("D10").Value = "25.06"
$ws.Range("E10").Value = "  +1.87%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.08395"
$ws.Range("E11").Value = "  +0.92%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.859.64"
$ws.Range("E12").Value = "  -0.23%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "5.231"
$ws.Range("E13").Value = "  +1.08%  "

# Row 14 - Polygon
$ws.Range("D14").Value = "0.7111"
$ws.Range("E14").Value = "  +0.39%  "

# Row 15 - Litecoin
$ws.Range("E15").Value = "  +0.08%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "29.320.73"
$ws.Range("E16").Value = "  +0.40%  "

# Row 17 - Uniswap
$ws.Range("D17").Value = "6.083"
$ws.Range("E17").Value = "  +2.95%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "0.000008220"
$ws.Range("E18").Value = "  +5.49%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "240.34"
$ws.Range("E19").Value = "  -0.84%  "

# Row 20 - Avalanche
$ws.Range("D20").Value = "13.19"

# Row 21 - WrappedliquidstakedEther2.0
$ws.Range("D21").Value = "2.121.77"
$ws.Range("E21").Value = "  -0.13%  "

# Row 22 - Dai
$ws.Range("D22").Value = "1.000"

# Row 23 - Chainlink
$ws.Range("D23").Value = "7.765"
$ws.Range("E23").Value = "  -1.03%  "

# Row 24 - BinanceUSD
$ws.Range("E24").Value = "  +0.01%  "

# Row 25 - Stellar
$ws.Range("D25").Value = "0.1600"
$ws.Range("E25").Value = "  +0.99%  "

# Row 26 - Monero
$ws.Range("D26").Value = "162.85"
$ws.Range("E26").Value = "  -0.14%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "9.024"
$ws.Range("E27").Value = "  +1.09%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "18.49"
$ws.Range("E28").Value = "  +0.40%  "

# Row 29 - PancakeSwap
$ws.Range("D29").Value = "1.508"
$ws.Range("E29").Value = "  +0.71%  "

# Row 30 - Filecoin
$ws.Range("D30").Value = "4.413"
$ws.Range("E30").Value = "  +0.46%  "

# Row 31 - Toncoin
$ws.Range("E31").Value = "  -2.66%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").Value = "4.306"

# Row 33 - Hedera
$ws.Range("D33").Value = "0.05294"
$ws.Range("E33").Value = "  +3.31%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  +1.66%  "

# Row 35 - ARBITRUM
$ws.Range("D35").Value = "1.175"
$ws.Range("E35").Value = "  +1.34%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "0.7446"
$ws.Range("E36").Value = "  -6.23%  "

# Row 37 - HuobiToken
$ws.Range("D37").Value = "2.700"
$ws.Range("E37").Value = "  +0.69%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "0.01868"
$ws.Range("E38").Value = "  +1.68%  "

# Row 39 - Maker
$ws.Range("D39").Value = "1.227.34"
$ws.Range("E39").Value = "  +5.32%  "

# Row 40 - MXToken
$ws.Range("D40").Value = "2.725"
$ws.Range("E40").Value = "  +0.98%  "

# Row 41 - FraxShare
$ws.Range("D41").Value = "6.564"
$ws.Range("E41").Value = "  +6.39%  "

# Row 42 - Quant
$ws.Range("D42").Value = "110.67"
$ws.Range("E42").Value = "  +8.40%  "

# Row 43 - TrustWalletToken
$ws.Range("D43").Value = "0.8856"
$ws.Range("E43").Value = "  -0.18%  "

# Row 44 - Aave
$ws.Range("D44").Value = "72.61"
$ws.Range("E44").Value = "  -0.26%  "

# Row 45 - PaxDollar
$ws.Range("D45").Value = "0.9998"
$ws.Range("E45").Value = "  -0.04%  "

# Row 46 - RocketPoolETH
$ws.Range("D46").Value = "2.018.79"
$ws.Range("E46").Value = "  +0.36%  "

# Row 47 - RenderToken
$ws.Range("D47").Value = "1.800"
$ws.Range("E47").Value = "  +2.09%  "

# Row 48 - Mantle
$ws.Range("D48").Value = "0.5197"
$ws.Range("E48").Value = "  +0.36%  "

# Row 49 - BabyDogeCoin
$ws.Range("E49").Value = "  +2.10%  "

# Row 50 - EnergySwap
$ws.Range("D50").Value = "9.397"
$ws.Range("E50").Value = "  +1.28%  "

# Row 51 - TheSandbox
$ws.Range("E51").Value = "  +1.32%  "

# Restore plain (unstyled) look now that the text values are locked in.
$ws.Range("D2:E51").Style = "Normal"
